$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.643996357917786
$ws.Range("B1").Value = 2.569762468338013
$ws.Range("C1").Value = 2.986007452011108
$ws.Range("D1").Value = 3.031073093414307
$ws.Range("E1").Value = 1.053208231925964
